$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-10-07 Monday" "2024-10-08 Tuesday"
Replace-Text "35-4=31" "31-18=13"
Replace-Text "49-32=17" "68-52=16"
Replace-Text "6+33=39" "27+24=51"
Replace-Text "62+37=99" "32-11=21"
Replace-Text "87-33=54" "0+59=59"
Replace-Text "30-21=9" "86-57=29"
Replace-Text "47+51=98" "85-26=59"
Replace-Text "97-38=59" "83-4=79"
Replace-Text "32+16=48" "74-70=4"
Replace-Text "22+76=98" "28-23=5"
Replace-Text "4+72=76" "85-60=25"
Replace-Text "98-43=55" "14+15=29"
Replace-Text "44-3=41" "59+33=92"
Replace-Text "2+22=24" "73-53=20"
Replace-Text "56-2=54" "97-35=62"
Replace-Text "54-50=4" "51+0=51"
Replace-Text "45+51=96" "80+0=80"
Replace-Text "31+35=66" "2+36=38"
Replace-Text "5+36=41" "58+2=60"
Replace-Text "25+44=69" "91+1=92"
Replace-Text "36+37=73" "90-59=31"
Replace-Text "87-43=44" "58-26=32"
Replace-Text "35+59=94" "89-80=9"
Replace-Text "25+17=42" "75-10=65"
Replace-Text "7+25=32" "7+37=44"
Replace-Text "57+17=74" "98-11=87"
Replace-Text "64+15=79" "69-33=36"
Replace-Text "20+6=26" "78-57=21"
Replace-Text "79-15=64" "6+18=24"
Replace-Text "32+48=80" "73+19=92"
Replace-Text "58-9=49" "24+58=82"
Replace-Text "84+10=94" "70+2=72"
Replace-Text "81-40=41" "14+21=35"
Replace-Text "40-10=30" "85+0=85"
Replace-Text "6+25=31" "17-16=1"
Replace-Text "41+44=85" "40+38=78"
Replace-Text "98-85=13" "91-10=81"
Replace-Text "25-16=9" "14+85=99"
Replace-Text "63-26=37" "55+41=96"
Replace-Text "92-67=25" "44+3=47"
Replace-Text "54-10=44" "30+63=93"
Replace-Text "51+33=84" "25-9=16"
Replace-Text "8+76=84" "79-39=40"
Replace-Text "65+1=66" "53+41=94"
Replace-Text "53+16=69" "29+48=77"
Replace-Text "6+44=50" "28-1=27"
Replace-Text "28-22=6" "13+18=31"
Replace-Text "38+34=72" "98-11=87"
Replace-Text "57-52=5" "3+7=10"
Replace-Text "86-51=35" "1+34=35"
Replace-Text "96-48=48" "72-41=31"
Replace-Text "81-44=37" "18+81=99"
Replace-Text "22+72=94" "75-15=60"
Replace-Text "77+4=81" "37+10=47"
Replace-Text "85-39=46" "83+9=92"
Replace-Text "96+2=98" "25+53=78"
Replace-Text "92-49=43" "91-84=7"
Replace-Text "78+1=79" "2+25=27"
Replace-Text "10+71=81" "89-79=10"
Replace-Text "5+31=36" "80-41=39"
Replace-Text "56+32=88" "36+19=55"
Replace-Text "13+33=46" "47-22=25"
Replace-Text "20+32=52" "52-40=12"
Replace-Text "81-12=69" "61+24=85"
Replace-Text "0+1=1" "97-10=87"
Replace-Text "19+15=34" "88-83=5"
Replace-Text "39+15=54" "16+81=97"
Replace-Text "38+40=78" "61-53=8"
Replace-Text "72+3=75" "12+22=34"
Replace-Text "19+3=22" "77-21=56"
Replace-Text "1+85=86" "53+28=81"
Replace-Text "78+14=92" "80+8=88"
Replace-Text "47-30=17" "54-20=34"
Replace-Text "23+24=47" "83-8=75"
Replace-Text "77-26=51" "71-56=15"
Replace-Text "44+5=49" "66+12=78"
Replace-Text "26-6=20" "46+24=70"
Replace-Text "54-12=42" "42+9=51"
Replace-Text "97-89=8" "52+31=83"
Replace-Text "84-48=36" "36+39=75"
Replace-Text "59-41=18" "86-43=43"
Replace-Text "30-12=18" "33+37=70"
Replace-Text "4+16=20" "64-4=60"
Replace-Text "9+56=65" "7+10=17"
Replace-Text "73-23=50" "58+21=79"
Replace-Text "38-15=23" "91-29=62"
Replace-Text "84-45=39" "6+51=57"
Replace-Text "80-23=57" "2+2=4"
Replace-Text "15+35=50" "34+49=83"
Replace-Text "18+44=62" "21+37=58"
Replace-Text "28+41=69" "69-49=20"
Replace-Text "73-25=48" "15+38=53"
Replace-Text "69-67=2" "20+44=64"
Replace-Text "78-55=23" "74-47=27"
Replace-Text "44-9=35" "99-65=34"
Replace-Text "61+6=67" "91-0=91"
Replace-Text "95-60=35" "77-38=39"
Replace-Text "16+49=65" "99-7=92"
Replace-Text "67-59=8" "0+94=94"
Replace-Text "12+71=83" "75+18=93"
